# Atualização de bases das ligas, do dia: 29-05-2024 às 22:54
#
# This edit permutes the betting-odds data (columns B:AD) among several
# rows of the single worksheet, while each row keeps its own "A" (id)
# column fixed. Concretely, row N ends up holding the B:AD content that
# originally belonged to a different row, as described by the following
# cycles (read as "row -> takes the B:AD content that used to be in row"):
#
#   85  <- 86,  86  <- 85                      (swap)
#   179 <- 180, 180 <- 179                     (swap)
#   232 <- 233, 233 <- 238, 238 <- 235, 235 <- 234,
#   234 <- 237, 237 <- 236, 236 <- 232          (7-cycle)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    return $ws.Range("B$row`:AD$row").Value()
}

# capture the "before" B:AD snapshot of every row that participates in a
# permutation, before any writes happen
$rows = @(85, 86, 179, 180, 232, 233, 234, 235, 236, 237, 238)
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = Get-RowData $r
}

# row (key) <- source row (value) whose original B:AD content it receives
$mapping = @{
    85  = 86
    86  = 85
    179 = 180
    180 = 179
    232 = 233
    233 = 238
    234 = 237
    235 = 234
    236 = 232
    237 = 236
    238 = 235
}

foreach ($r in $rows) {
    $src = $mapping[$r]
    $ws.Range("B$r`:AD$r").Value = $snapshot[$src]
}
